$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ref, $val)
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '26.982.36'
$ws.Range('D3').Value = '1.676.48'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-CellText D5 '215.10'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +0.29%  '
Set-CellText D10 '20.32'
$ws.Range('E10').Value = '  +0.72%  '
Set-CellText D11 '0.0887'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('D12').Value = '1.912.58'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '1.686.04'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('E15').Value = '  +1.54%  '
Set-CellText D16 '65.81'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '26.995.34'
$ws.Range('E17').Value = '  -0.24%  '
Set-CellText D18 '237.12'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('E19').Value = '  +5.34%  '
$ws.Range('D20').Value = '0.0₃0733'
$ws.Range('E20').Value = '  -0.70%  '
Set-CellText D22 '4.44'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -1.75%  '
Set-CellText D25 '145.99'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('E27').Value = '  +1.56%  '
Set-CellText D28 '0.113'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  +0.07%  '
Set-CellText D30 '0.0498'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -0.39%  '
Set-CellText D32 '3.33'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '1.479.55'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('E36').Value = '  +0.05%  '
Set-CellText D37 '0.585'
$ws.Range('E37').Value = '  +2.34%  '
$ws.Range('E38').Value = '  +2.21%  '
Set-CellText D39 '0.904'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  +1.81%  '
Set-CellText D44 '67.47'
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').Value = '1.819.13'
Set-CellText D46 '0.782'
$ws.Range('E46').Value = '  +0.36%  '
Set-CellText D47 '90.49'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('E50').Value = '  +1.67%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText D51 '0.0508'
$ws.Range('E51').Value = '  -0.03%  '
